$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Small fish fry* {1.5}'
$ws.Range('E2').Value = 'Small fish fry*'
$ws.Range('D3').Value = 'Ruti* {1.0}; Onion, raw {1.8}'
$ws.Range('D4').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {0.53}; Pangas, without bones, raw {0.53}; Potato, Diamond, boiled* (without salt) {0.97}; Brinjal, purple, long, boiled* (without salt) {0.97}'
$ws.Range('B5').Value = 'chicken (ootato with bread)'
$ws.Range('D5').Value = 'Chicken leg, without skin, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.97}; Bread, bun/roll {0.97}'
$ws.Range('E5').Value = 'Chicken leg, without skin, raw'
$ws.Range('D7').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Egg, chicken, farmed, boiled* (without salt) {0.17}; Potato, Diamond, boiled* (without salt) {0.16}; Onion, raw {0.16}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Soybean oil {0.15}'
$ws.Range('D8').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Egg, chicken, farmed, boiled* (without salt) {1.05}; Potato, Diamond, boiled* (without salt) {0.97}; Onion, raw {0.97}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Soybean oil {0.15}'
$ws.Range('D10').Value = 'UNKNOWN; UNKNOWN; Potato, Diamond, boiled* (without salt) {0.9}; Papaya, unripe, boiled* (without salt) {0.9}; Green gram, split, boiled* (without salt) {0.9}; Soybean oil {0.15}; Onion, raw {0.9}; Garlic, raw {0.9}'
$ws.Range('D12').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Egg, chicken, farmed, boiled* (without salt) {0.17}; Onion, raw {0.33}; Soybean oil {0.15}'
$ws.Range('D14').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Egg, chicken, farmed, boiled* (without salt) {0.5}; Onion, raw {1.12}; Soybean oil {0.15}; Green gram, split, boiled* (without salt) {1.12}; Gourd, ash, raw {1.12}; Prawn, Giant tiger prawn, raw {1.12}'
$ws.Range('D16').Value = 'Bread, bun/roll {1.0}; Egg, chicken, farmed, boiled* (without salt) {0.17}; Tomato, red, ripe, boiled* (without salt) {0.11}; Chilli, green, with seeds, raw {0.11}; Onion, raw {0.11}; Turmeric, dried {0.1}; Soybean oil {0.15}; Water, drinking {2.5}'
$ws.Range('D18').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Chicken leg, without skin, raw {1.57}; Potato, Diamond, boiled* (without salt) {2.92}; Water, drinking {2.5}'
$ws.Range('D19').Value = 'Chicken leg, without skin, raw {1.57}; Potato, Diamond, boiled* (without salt) {2.92}; Banana, Sagar, ripe, raw {0.5}; Biscuit, sweet* {0.5}; Orange juice, raw (unsweetened) {3.0}'
$ws.Range('D21').Value = 'Bread, bun/roll {1.0}; Potato, Diamond, boiled* (without salt) {9.0}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Soybean oil {0.15}; Water, drinking {5.0}'
$ws.Range('D23').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Chicken leg, without skin, raw {9.0}; Water, drinking {2.5}'
$ws.Range('D25').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Chicken leg, without skin, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.97}; Green gram, split, boiled* (without salt) {0.97}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Salt {0.1}'
$ws.Range('D27').Value = 'Bread, bun/roll {1.0}; Potato, Diamond, boiled* (without salt) {1.5}; Chilli, red, dry {0.1}; Turmeric, dried {0.1}; Onion, raw {1.5}; Soybean oil {0.15}; Mustard oil {0.15}; Water, drinking {2.5}'
$ws.Range('D28').Value = 'UNKNOWN; Emblic, raw {0.5}; Water, drinking {1.25}'
$ws.Range('D29').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {1.57}; Potato, Diamond, boiled* (without salt) {0.73}; Onion, raw {0.73}; Chilli, red, dry {0.1}; Soybean, dried, raw {0.73}; Soybean oil {0.15}; Cauliflower, boiled* (without salt) {0.73}'
$ws.Range('D33').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.97}; Onion, raw {0.97}; Chilli, red, dry {0.1}; Turmeric, dried {0.1}; Water, drinking {2.5}'
$ws.Range('E33').Value = 'Pangas, without bones, raw'
$ws.Range('D35').Value = 'Bread, bun/roll {1.0}; Potato, Diamond, boiled* (without salt) {1.5}; Tomato, red, ripe, boiled* (without salt) {1.5}; Soybean oil {0.15}; Chilli, green, with seeds, raw {1.5}; Coriander seed, dry {0.1}; Coriander leaves, raw {0.1}; Water, drinking {2.5}'
$ws.Range('D36').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {0.17}; Amaranth, leaves, red, boiled* (without salt) {0.33}; Coriander leaves, raw {0.1}'
$ws.Range('D40').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Potato, Diamond, boiled* (without salt) {0.5}; Chilli, red, dry {0.1}; Onion, raw {0.5}; Soybean oil {0.15}; Green gram, split, boiled* (without salt) {0.5}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Onion, raw {0.5}; Salt {0.1}; Orange juice, raw (unsweetened) {5.0}'
$ws.Range('D42').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Potato, Diamond, boiled* (without salt) {1.38}; Chilli, red, dry {0.1}; Onion, raw {1.38}; Soybean oil {0.15}; Green gram, split, boiled* (without salt) {1.38}; Turmeric, dried {0.1}; Chilli, red, dry {0.1}; Onion, raw {1.38}; Salt {0.1}; Water, drinking {5.0}'
$ws.Range('D46').Value = 'Bread, bun/roll {1.0}; Potato, Diamond, boiled* (without salt) {2.25}; Onion, raw {2.25}; Soybean oil {0.15}; Water, drinking {5.0}'
$ws.Range('D48').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {3.15}; Tomato, red, ripe, boiled* (without salt) {1.17}; Onion, raw {1.17}; Soybean oil {0.15}; Potato, Diamond, boiled* (without salt) {1.17}; Onion, raw {1.17}; Soybean oil {0.15}; Pumpkin leaves, raw {1.17}; Water, drinking {2.5}'
$ws.Range('E48').Value = 'Pangas, without bones, raw'
$ws.Range('D49').Value = 'Potato, Diamond, boiled* (without salt) {0.43}; Carrot, boiled* (without salt) {0.43}; Onion, raw {0.43}; Soybean oil {0.15}; Chilli, green, with seeds, raw {0.43}; Potato, Diamond, boiled* (without salt) {0.43}; Onion, raw {0.43}; Chilli, green, with seeds, raw {0.43}; Soybean oil {0.15}; UNKNOWN; Biscuit, sweet* {0.75}'
$ws.Range('D50').Value = 'Egg, chicken, farmed, boiled* (without salt) {0.17}; Onion, raw {0.11}; Chilli, red, dry {0.1}; Coriander leaves, raw {0.1}; Pangas, without bones, raw {0.17}; Soybean oil {0.15}; Onion, raw {0.11}; Tomato, red, ripe, boiled* (without salt) {0.11}; Green gram, split, boiled* (without salt) {0.11}; Onion, raw {0.11}; Coriander leaves, raw {0.1}; Garlic, raw {0.11}; Soybean oil {0.15}; Rice, BR-28, boiled* (without salt) {1.5}; Water, drinking {1.25}'
$ws.Range('E50').Value = 'Pangas, without bones, raw'
$ws.Range('D52').Value = 'Bread, bun/roll {1.0}; Egg, chicken, farmed, boiled* (without salt) {1.57}; Tomato, red, ripe, boiled* (without salt) {0.97}; Chilli, green, with seeds, raw {0.97}; Onion, raw {0.97}; Turmeric, dried {0.1}; Soybean oil {0.15}; Water, drinking {2.5}'
$ws.Range('D53').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.65}; Onion, raw {0.65}; Chilli, red, dry {0.1}; Soybean oil {0.15}; Cauliflower, boiled* (without salt) {0.65}; UNKNOWN; Water, drinking {5.0}'
$ws.Range('E53').Value = 'Pangas, without bones, raw'
$ws.Range('D55').Value = 'Biscuit, sweet* {4.5}'
$ws.Range('D56').Value = 'Pangas, without bones, raw {1.05}; Tomato, red, ripe, boiled* (without salt) {0.65}; Radish, boiled* (without salt) {0.65}; UNKNOWN; Cucumber, peeled, raw {0.65}; Soybean oil {0.15}; UNKNOWN'
$ws.Range('E56').Value = 'Pangas, without bones, raw'
$ws.Range('D59').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.65}; Turmeric, dried {0.1}; Chilli, green, with seeds, raw {0.65}; Soybean oil {0.15}; Green gram, split, boiled* (without salt) {0.65}; Chilli, red, dry {0.1}; Turmeric, dried {0.1}'
$ws.Range('E59').Value = 'Pangas, without bones, raw'
$ws.Range('D60').Value = 'Milk, cow, whole fat (pasteurized, UTH) {0.75}; Rice, BR-28, boiled* (without salt) {0.75}; Jaggery, sugarcane, solid {0.75}; Payesh* {0.75}'
$ws.Range('D61').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Pangas, without bones, raw {1.05}; Potato, Diamond, boiled* (without salt) {0.33}; Soybean oil {0.15}; Green gram, split, boiled* (without salt) {0.33}; Chilli, red, dry {0.1}; Turmeric, dried {0.1}; Brinjal, purple, long, boiled* (without salt) {0.33}; Potato, Diamond, boiled* (without salt) {0.33}; Radish, boiled* (without salt) {0.33}; Spinach, boiled* (without salt) {0.33}'
$ws.Range('E61').Value = 'Pangas, without bones, raw'
$ws.Range('D63').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Spinach, boiled* (without salt) {4.5}; Soybean oil {0.15}; Water, drinking {2.5}'
$ws.Range('D64').Value = 'Rice, BR-28, boiled* (without salt) {3.0}; Chicken leg, without skin, raw {1.57}; Onion, raw {1.46}; Garlic, raw {1.46}; Ginger root, raw {0.1}; Chilli, red, dry {0.1}; Water, drinking {2.5}'
$ws.Range('E64').Value = 'Chicken leg, without skin, raw'
$ws.Range('D65').Value = 'UNKNOWN; Orange juice, raw (unsweetened) {1.25}'
$ws.Range('D66').Value = 'Hog plum, raw {3.0}; Salt {0.1}; Chilli, red, dry {0.1}'
$ws.Range('D67').Value = 'Rice, puffed, salted {2.25}; Tomato, red, ripe, boiled* (without salt) {2.25}; Chilli, red, dry {0.1}; Coriander leaves, raw {0.1}; Water, drinking {5.0}'
